$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column O: "Identifiers" -> "Internal  House  Name" header, and
#     the "Internal House Name:" prefix dropped from the data rows
#     (the prefix now lives in the header instead) ---
$ws.Range("O1").Value = "Internal  House  Name"
$ws.Range("O2").Value = "a round nose"
$ws.Range("O4").Value = "a square nose"
$ws.Range("O5").Value = "a yellow nose"

# --- New column S: "Tag /Band" header (tracks a split-out identifiers column) ---
$ws.Range("S1").Value = "Tag /Band"

# Give column O a bit more breathing room now that the header text is longer,
# and put the viewport/selection roughly where the editor left it.
$ws.Range("O1:O5").ColumnWidth = 14.6
[void]$ws.Range("O1").Select()

Write-Output "done"
